$d = $word.ActiveDocument

# Locate the single paragraph that documents TankShooting's responsibility.
$rng = $d.Content
$found = $rng.Find.Execute("TankShooting is responsible for load the subtank")
if (-not $found) {
    throw "Could not find target paragraph text"
}

# Collapse to the start of the match: InsertXML-ing a <w:p>-wrapped fragment at
# any point inside a paragraph swaps that whole paragraph (text + bookmark)
# out for the inserted paragraph(s), reusing the paragraph mark that follows.
$rng.Collapse(1)

# Split the old single paragraph into:
#   1) a new paragraph explaining why SetDynamicObjectLibrary/OnChangeTank live apart
#   2) a blank paragraph
#   3) the original sentence, now prefixed with "TankShooting(100)" (bookmark kept
#      in place) and suffixed with a note about loading time
$xml = @'
<?xml version="1.0" encoding="utf-8"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>SetDynamicObjectLibrary</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>和</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>OnChangeTank</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>不在一起，是因为</w:t></w:r><w:r><w:t>SetDynamicObjectLibrary</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>只用</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>load</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>一次。</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:t>TankShooting</w:t></w:r><w:r><w:t>(100)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> is responsible for load the subtank. It is the most important one among TankMovement, Tankhealth, Tankshooting.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>注意</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">loading time. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($xml)
Write-Output "Inserted. Paragraph count: $($d.Paragraphs.Count)"
